$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '49.710.34'
$ws.Range('E2').Value = '  +2.92%  '

$ws.Range('D3').Value = '2.624.72'
$ws.Range('E3').Value = '  +4.71%  '

$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').Value = '325.50'
$ws.Range('E5').Value = '  +1.30%  '

$ws.Range('D6').Value = '109.85'
$ws.Range('E6').Value = '  +1.48%  '

$ws.Range('D7').Value = '0.534'
$ws.Range('E7').Value = '  +1.04%  '

$ws.Range('E8').Value = '  -0.11%  '

$ws.Range('D9').Value = '0.560'
$ws.Range('E9').Value = '  +3.39%  '

$ws.Range('D10').Value = '40.65'
$ws.Range('E10').Value = '  +1.78%  '

$ws.Range('D11').Value = '20.70'
$ws.Range('E11').Value = '  +2.11%  '

$ws.Range('D12').Value = '0.0823'
$ws.Range('E12').Value = '  +0.46%  '

$ws.Range('E13').Value = '  +0.69%  '

$ws.Range('D14').Value = '7.29'
$ws.Range('E14').Value = '  +1.41%  '

$ws.Range('D15').Value = '3.034.08'
$ws.Range('E15').Value = '  +4.71%  '

$ws.Range('D16').Value = '2.619.65'
$ws.Range('E16').Value = '  +4.65%  '

$ws.Range('D17').Value = '0.873'
$ws.Range('E17').Value = '  +3.05%  '

$ws.Range('D18').Value = '49.663.20'
$ws.Range('E18').Value = '  +3.13%  '

$ws.Range('D19').Value = '3.11'
$ws.Range('E19').Value = '  +11.78%  '

$ws.Range('D20').Value = '13.36'
$ws.Range('E20').Value = '  +1.85%  '

$ws.Range('D21').Value = '6.80'
$ws.Range('E21').Value = '  +0.69%  '

$ws.Range('D22').Value = '0.0₃0953'
$ws.Range('E22').Value = '  +0.63%  '

$ws.Range('D23').Value = '281.21'
$ws.Range('E23').Value = '  +0.17%  '

$ws.Range('D24').Value = '72.74'
$ws.Range('E24').Value = '  +0.72%  '

$ws.Range('D25').Value = '2.59'
$ws.Range('E25').Value = '  +1.24%  '

$ws.Range('D26').Value = '26.65'
$ws.Range('E26').Value = '  +3.39%  '

$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.05%  '

$ws.Range('D28').Value = '2.24'
$ws.Range('E28').Value = '  -6.71%  '

$ws.Range('D29').Value = '10.00'
$ws.Range('E29').Value = '  +2.00%  '

$ws.Range('D30').Value = '0.145'
$ws.Range('E30').Value = '  +2.56%  '

$ws.Range('D31').Value = '36.19'
$ws.Range('E31').Value = '  +2.38%  '

$ws.Range('D32').Value = '49.65'
$ws.Range('E32').Value = '  +0.61%  '

$ws.Range('D33').Value = '19.69'
$ws.Range('E33').Value = '  +0.44%  '

$ws.Range('D34').Value = '5.45'
$ws.Range('E34').Value = '  +1.59%  '

$ws.Range('E35').Value = '  -0.09%  '

$ws.Range('D36').Value = '0.0794'
$ws.Range('E36').Value = '  +1.22%  '

$ws.Range('D37').Value = '2.06'
$ws.Range('E37').Value = '  +5.43%  '

$ws.Range('D38').Value = '4.74'
$ws.Range('E38').Value = '  +1.77%  '

$ws.Range('E39').Value = '  +5.86%  '

$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').Value = '22.76'
$ws.Range('E40').Value = '  +5.24%  '

$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').Value = '0.113'
$ws.Range('E41').Value = '  +0.74%  '

$ws.Range('D42').Value = '123.40'
$ws.Range('E42').Value = '  +1.73%  '

$ws.Range('D43').Value = '2.22'
$ws.Range('E43').Value = '  +0.27%  '

$ws.Range('E44').Value = '  +3.48%  '

$ws.Range('D45').Value = '3.37'
$ws.Range('E45').Value = '  +6.09%  '

$ws.Range('D46').Value = '2.056.56'
$ws.Range('E46').Value = '  +2.10%  '

$ws.Range('D47').Value = '2.23'
$ws.Range('E47').Value = '  +11.92%  '

$ws.Range('E48').Value = '  +8.89%  '

$ws.Range('D49').Value = '9.04'
$ws.Range('E49').Value = '  +0.68%  '

$ws.Range('D50').Value = '5.38'
$ws.Range('E50').Value = '  +3.76%  '

$ws.Range('D51').Value = '81.93'
$ws.Range('E51').Value = '  +1.77%  '
